# Turn off automatic hyphenation (adds <w:suppressAutoHyphens w:val="true"/>
# to <w:pPr>) on the built-in paragraph styles: Normal, Heading 1-6 and
# LO-normal. This mirrors Word's Format > Paragraph > Line and Page Breaks >
# "Don't hyphenate" checkbox applied at the style level.

$d = $word.ActiveDocument

$styleNames = @(
    "Normal",
    "Heading 1",
    "Heading 2",
    "Heading 3",
    "Heading 4",
    "Heading 5",
    "Heading 6",
    "LO-normal"
)

foreach ($styleName in $styleNames) {
    $style = $d.Styles($styleName)
    $style.ParagraphFormat.Hyphenation = $false
}
